# Refresh the cryptos list with newly scraped price / 1h-volume figures.
# (Mirrors the nightly GitHub Actions scraper commit.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the Price column stores plain-looking numeric strings ("230.67",
# "0.619", ...) as TEXT (matches the source scrape, which also renders
# multi-dot thousands values like "35.050.05" as text). Assigning a bare
# numeric-looking string to .Value auto-coerces it to a real number, so we
# force the cell to Text format first for any value Excel would otherwise
# treat as numeric. Multi-dot strings are already unambiguous text and
# don't need this.

$ws.Range("D2").Value = "35.050.05"
$ws.Range("E2").Value = "  +0.34%  "

$ws.Range("D3").Value = "1.820.00"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.67"
$ws.Range("E5").Value = "  -0.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  +0.87%  "

$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.84"
$ws.Range("E8").Value = "  -4.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.323"
$ws.Range("E9").Value = "  +4.78%  "

$ws.Range("E10").Value = "  +0.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0990"
$ws.Range("E11").Value = "  -1.18%  "

$ws.Range("D12").Value = "2.083.49"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("E13").Value = "  +2.62%  "

$ws.Range("E14").Value = "  +1.44%  "

$ws.Range("D15").Value = "1.819.59"
$ws.Range("E15").Value = "  -0.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.65"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").Value = "35.045.48"
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.70"
$ws.Range("E18").Value = "  +0.39%  "

$ws.Range("E19").Value = "  +0.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.06"
$ws.Range("E20").Value = "  +1.25%  "

$ws.Range("E21").Value = "  +2.80%  "

$ws.Range("E22").Value = "  +2.76%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("E24").Value = "  +0.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "173.51"
$ws.Range("E25").Value = "  +0.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.82"
$ws.Range("E26").Value = "  +0.96%  "

$ws.Range("E27").Value = "  +3.53%  "

$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.52"
$ws.Range("E29").Value = "  -2.61%  "

$ws.Range("E30").Value = "  -0.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.00"
$ws.Range("E31").Value = "  +3.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0551"
$ws.Range("E32").Value = "  +0.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.97"
$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.26"
$ws.Range("E34").Value = "  +13.59%  "

$ws.Range("E35").Value = "  +3.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.696"
$ws.Range("E36").Value = "  +3.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "92.76"
$ws.Range("E37").Value = "  +0.60%  "

$ws.Range("E38").Value = "  +6.46%  "

$ws.Range("D39").Value = "1.341.03"
$ws.Range("E39").Value = "  +2.10%  "

$ws.Range("E40").Value = "  +1.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.985"
$ws.Range("E41").Value = "  +0.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.75"
$ws.Range("E42").Value = "  +0.46%  "

$ws.Range("E43").Value = "  -0.33%  "

$ws.Range("E44").Value = "  -1.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.76"
$ws.Range("E45").Value = "  -0.44%  "

# Row 46 (was Kaspa) becomes FraxShare
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.23"
$ws.Range("E46").Value = "  +0.93%  "

# Row 47 (was FraxShare) becomes Kaspa
$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0519"
$ws.Range("E47").Value = "  +2.05%  "

$ws.Range("D48").Value = "2.000.24"
$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0666"
$ws.Range("E50").Value = "  +4.36%  "

# Row 51 (was Quant) becomes THORChain
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.21"
$ws.Range("E51").Value = "  +13.64%  "
